# Apply updated coin data (price/volume/symbol refresh) per commit
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("D2:E2")
$rng.NumberFormat = "@"
$ws.Range("D2").Value = "244.99"
$ws.Range("E2").Value = "-0.65%"
$rng.ClearFormats()

$rng = $ws.Range("D3:E3")
$rng.NumberFormat = "@"
$ws.Range("D3").Value = "26.94"
$ws.Range("E3").Value = "1.50%"
$rng.ClearFormats()

$rng = $ws.Range("D4:E4")
$rng.NumberFormat = "@"
$ws.Range("D4").Value = "5.075"
$ws.Range("E4").Value = "-0.27%"
$rng.ClearFormats()

$rng = $ws.Range("D5:E5")
$rng.NumberFormat = "@"
$ws.Range("D5").Value = "0.05691"
$ws.Range("E5").Value = "1.41%"
$rng.ClearFormats()

$rng = $ws.Range("D6:E6")
$rng.NumberFormat = "@"
$ws.Range("D6").Value = "6.479"
$ws.Range("E6").Value = "-0.05%"
$rng.ClearFormats()

$rng = $ws.Range("D7:E7")
$rng.NumberFormat = "@"
$ws.Range("D7").Value = "0.8202"
$ws.Range("E7").Value = "0.85%"
$rng.ClearFormats()

$rng = $ws.Range("D8:E8")
$rng.NumberFormat = "@"
$ws.Range("D8").Value = "0.8386"
$ws.Range("E8").Value = "-0.82%"
$rng.ClearFormats()

$rng = $ws.Range("D9:E9")
$rng.NumberFormat = "@"
$ws.Range("D9").Value = "0.1328"
$ws.Range("E9").Value = "-0.89%"
$rng.ClearFormats()

$rng = $ws.Range("D10:E10")
$rng.NumberFormat = "@"
$ws.Range("D10").Value = "0.06904"
$ws.Range("E10").Value = "-0.92%"
$rng.ClearFormats()

$rng = $ws.Range("D11:E11")
$rng.NumberFormat = "@"
$ws.Range("D11").Value = "0.02858"
$ws.Range("E11").Value = "0.20%"
$rng.ClearFormats()

$rng = $ws.Range("D12:E12")
$rng.NumberFormat = "@"
$ws.Range("D12").Value = "0.09402"
$ws.Range("E12").Value = "0.15%"
$rng.ClearFormats()

$rng = $ws.Range("D13:E13")
$rng.NumberFormat = "@"
$ws.Range("D13").Value = "0.001521"
$ws.Range("E13").Value = "0.20%"
$rng.ClearFormats()

$rng = $ws.Range("D14:E14")
$rng.NumberFormat = "@"
$ws.Range("D14").Value = "0.04102"
$ws.Range("E14").Value = "-12.11%"
$rng.ClearFormats()

$rng = $ws.Range("B15:E15")
$rng.NumberFormat = "@"
$ws.Range("B15").Value = "One"
$ws.Range("C15").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D15").Value = "0.0005993"
$ws.Range("E15").Value = "-93.93%"
$rng.ClearFormats()

$rng = $ws.Range("B16:E16")
$rng.NumberFormat = "@"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "0.006094"
$ws.Range("E16").Value = "-0.25%"
$rng.ClearFormats()

$rng = $ws.Range("B17:E17")
$rng.NumberFormat = "@"
$ws.Range("B17").Value = "UpBots"
$ws.Range("C17").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("D17").Value = "0.007486"
$ws.Range("E17").Value = "3,761.28%"
$rng.ClearFormats()

$rng = $ws.Range("B18:E18")
$rng.NumberFormat = "@"
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").Value = "3.510"
$ws.Range("E18").Value = "-2.25%"
$rng.ClearFormats()

$rng = $ws.Range("B19:E19")
$rng.NumberFormat = "@"
$ws.Range("B19").Value = "GateToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D19").Value = "3.002"
$ws.Range("E19").Value = "-0.29%"
$rng.ClearFormats()

$rng = $ws.Range("B20:E20")
$rng.NumberFormat = "@"
$ws.Range("B20").Value = "BTSEToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D20").Value = "2.316"
$ws.Range("E20").Value = "9.32%"
$rng.ClearFormats()

$rng = $ws.Range("B21:E21")
$rng.NumberFormat = "@"
$ws.Range("B21").Value = "BitpandaEcosystemToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D21").Value = "0.3177"
$ws.Range("E21").Value = "-0.21%"
$rng.ClearFormats()

$rng = $ws.Range("B22:E22")
$rng.NumberFormat = "@"
$ws.Range("B22").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C22").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D22").Value = "0.03175"
$ws.Range("E22").Value = "1.03%"
$rng.ClearFormats()

$rng = $ws.Range("B23:E23")
$rng.NumberFormat = "@"
$ws.Range("B23").Value = "ProBitToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D23").Value = "0.1297"
$ws.Range("E23").Value = "-1.78%"
$rng.ClearFormats()

$rng = $ws.Range("B24:E24")
$rng.NumberFormat = "@"
$ws.Range("B24").Value = "MCDex"
$ws.Range("C24").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D24").Value = "3.548"
$ws.Range("E24").Value = "-5.42%"
$rng.ClearFormats()

$rng = $ws.Range("B25:E25")
$rng.NumberFormat = "@"
$ws.Range("B25").Value = "ZBToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D25").Value = "0.1373"
$ws.Range("E25").Value = "1.64%"
$rng.ClearFormats()

$rng = $ws.Range("D26:E26")
$rng.NumberFormat = "@"
$ws.Range("D26").Value = "0.001218"
$ws.Range("E26").Value = "-2.61%"
$rng.ClearFormats()

$rng = $ws.Range("D27:E27")
$rng.NumberFormat = "@"
$ws.Range("D27").Value = "0.003967"
$ws.Range("E27").Value = "-13.55%"
$rng.ClearFormats()

$rng = $ws.Range("E28")
$rng.NumberFormat = "@"
$ws.Range("E28").Value = "1.97%"
$rng.ClearFormats()

$rng = $ws.Range("D40:E40")
$rng.NumberFormat = "@"
$ws.Range("D40").Value = "0.03696"
$ws.Range("E40").Value = "0.74%"
$rng.ClearFormats()

$rng = $ws.Range("D41:E41")
$rng.NumberFormat = "@"
$ws.Range("D41").Value = "0.005955"
$ws.Range("E41").Value = "-3.04%"
$rng.ClearFormats()

$rng = $ws.Range("D42:E42")
$rng.NumberFormat = "@"
$ws.Range("D42").Value = "0.1057"
$ws.Range("E42").Value = "0.07%"
$rng.ClearFormats()

$rng = $ws.Range("D43:E43")
$rng.NumberFormat = "@"
$ws.Range("D43").Value = "0.002298"
$ws.Range("E43").Value = "-11.75%"
$rng.ClearFormats()

$rng = $ws.Range("D44:E44")
$rng.NumberFormat = "@"
$ws.Range("D44").Value = "0.009392"
$ws.Range("E44").Value = "5.47%"
$rng.ClearFormats()

$rng = $ws.Range("D45:E45")
$rng.NumberFormat = "@"
$ws.Range("D45").Value = "0.00005190"
$ws.Range("E45").Value = "-1.98%"
$rng.ClearFormats()

$rng = $ws.Range("E46")
$rng.NumberFormat = "@"
$ws.Range("E46").Value = "-0.11%"
$rng.ClearFormats()

$rng = $ws.Range("E47")
$rng.NumberFormat = "@"
$ws.Range("E47").Value = "-15.53%"
$rng.ClearFormats()

$rng = $ws.Range("D48:E48")
$rng.NumberFormat = "@"
$ws.Range("D48").Value = "0.002594"
$ws.Range("E48").Value = "3.06%"
$rng.ClearFormats()

$rng = $ws.Range("E49")
$rng.NumberFormat = "@"
$ws.Range("E49").Value = "-0.11%"
$rng.ClearFormats()

$rng = $ws.Range("D50:E50")
$rng.NumberFormat = "@"
$ws.Range("D50").Value = "0.0001999"
$ws.Range("E50").Value = "-0.11%"
$rng.ClearFormats()
